$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInclude = $wb.Worksheets.Item("Include #0")

# URL value changed
$wsMeta.Range("B2").Value = "https://www.hl7.org/fhir/R4/valueset-allergyintolerance-verification"

# Title value changed
$wsMeta.Range("B5").Value = "NG-Imm AEFI Verification Status VS"

# Date value changed
$wsMeta.Range("B8").Value = "2025-06-25T06:29:04+01:00"

# System URI value now points to the same string as the URL above
$wsInclude.Range("B4").Value = "https://www.hl7.org/fhir/R4/valueset-allergyintolerance-verification"
